$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$dCell = $ws.Cells.Item(2, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '67.725.26'
$dCell.Style = $dStyle
$ws.Range('E2').Value = '  +2.36%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$dCell = $ws.Cells.Item(3, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '3.577.94'
$dCell.Style = $dStyle
$ws.Range('E3').Value = '  +1.38%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$dCell = $ws.Cells.Item(4, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '1.00'
$dCell.Style = $dStyle
$ws.Range('E4').Value = '  +0.30%  '

$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$dCell = $ws.Cells.Item(5, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '201.19'
$dCell.Style = $dStyle
$ws.Range('E5').Value = '  +8.26%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$dCell = $ws.Cells.Item(6, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '570.29'
$dCell.Style = $dStyle
$ws.Range('E6').Value = '  +0.23%  '

$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$dCell = $ws.Cells.Item(7, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '3.578.98'
$dCell.Style = $dStyle
$ws.Range('E7').Value = '  +1.55%  '

$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$dCell = $ws.Cells.Item(8, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.613'
$dCell.Style = $dStyle
$ws.Range('E8').Value = '  +1.17%  '

$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$dCell = $ws.Cells.Item(9, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '1.00'
$dCell.Style = $dStyle
$ws.Range('E9').Value = '  +0.13%  '

$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$dCell = $ws.Cells.Item(10, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.676'
$dCell.Style = $dStyle
$ws.Range('E10').Value = '  +2.52%  '

$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$dCell = $ws.Cells.Item(11, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '60.14'
$dCell.Style = $dStyle
$ws.Range('E11').Value = '  +13.49%  '

$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$dCell = $ws.Cells.Item(12, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.147'
$dCell.Style = $dStyle
$ws.Range('E12').Value = '  +1.96%  '

$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$dCell = $ws.Cells.Item(13, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.0000280'
$dCell.Style = $dStyle
$ws.Range('E13').Value = '  +9.77%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$dCell = $ws.Cells.Item(14, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '10.25'
$dCell.Style = $dStyle
$ws.Range('E14').Value = '  +6.05%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$dCell = $ws.Cells.Item(15, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '4.155.97'
$dCell.Style = $dStyle
$ws.Range('E15').Value = '  +1.88%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$dCell = $ws.Cells.Item(16, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '3.579.99'
$dCell.Style = $dStyle
$ws.Range('E16').Value = '  +1.98%  '

$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$dCell = $ws.Cells.Item(17, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.127'
$dCell.Style = $dStyle
$ws.Range('E17').Value = '  +1.25%  '

$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$dCell = $ws.Cells.Item(18, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '19.02'
$dCell.Style = $dStyle
$ws.Range('E18').Value = '  +5.37%  '

$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$dCell = $ws.Cells.Item(19, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '67.562.36'
$dCell.Style = $dStyle
$ws.Range('E19').Value = '  +2.40%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$dCell = $ws.Cells.Item(20, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '12.18'
$dCell.Style = $dStyle
$ws.Range('E20').Value = '  +1.99%  '

$ws.Range('B21').Value = 'Polygon'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$dCell = $ws.Cells.Item(21, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '1.06'
$dCell.Style = $dStyle
$ws.Range('E21').Value = '  +0.82%  '

$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$dCell = $ws.Cells.Item(22, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '402.26'
$dCell.Style = $dStyle
$ws.Range('E22').Value = '  +4.44%  '

$ws.Range('B23').Value = 'RenderToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$dCell = $ws.Cells.Item(23, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '12.69'
$dCell.Style = $dStyle
$ws.Range('E23').Value = '  +16.16%  '

$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$dCell = $ws.Cells.Item(24, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '4.19'
$dCell.Style = $dStyle
$ws.Range('E24').Value = '  +0.37%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$dCell = $ws.Cells.Item(25, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '84.51'
$dCell.Style = $dStyle
$ws.Range('E25').Value = '  -0.11%  '

$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$dCell = $ws.Cells.Item(26, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '2.88'
$dCell.Style = $dStyle
$ws.Range('E26').Value = '  +0.36%  '

$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$dCell = $ws.Cells.Item(27, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '12.41'
$dCell.Style = $dStyle
$ws.Range('E27').Value = '  +1.93%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$dCell = $ws.Cells.Item(28, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '3.85'
$dCell.Style = $dStyle
$ws.Range('E28').Value = '  +9.62%  '

$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$dCell = $ws.Cells.Item(29, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '9.17'
$dCell.Style = $dStyle
$ws.Range('E29').Value = '  +4.51%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$dCell = $ws.Cells.Item(30, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '7.63'
$dCell.Style = $dStyle
$ws.Range('E30').Value = '  +4.04%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$dCell = $ws.Cells.Item(31, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '31.36'
$dCell.Style = $dStyle
$ws.Range('E31').Value = '  +1.90%  '

$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$dCell = $ws.Cells.Item(32, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '677.71'
$dCell.Style = $dStyle
$ws.Range('E32').Value = '  +10.27%  '

$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$dCell = $ws.Cells.Item(33, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '12.10'
$dCell.Style = $dStyle
$ws.Range('E33').Value = '  -0.27%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$dCell = $ws.Cells.Item(34, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.113'
$dCell.Style = $dStyle
$ws.Range('E34').Value = '  +0.35%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$dCell = $ws.Cells.Item(35, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '63.26'
$dCell.Style = $dStyle
$ws.Range('E35').Value = '  +0.46%  '

$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$dCell = $ws.Cells.Item(36, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '41.37'
$dCell.Style = $dStyle
$ws.Range('E36').Value = '  +0.05%  '

$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$dCell = $ws.Cells.Item(37, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.408'
$dCell.Style = $dStyle
$ws.Range('E37').Value = '  +2.65%  '

$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$dCell = $ws.Cells.Item(38, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '1.00'
$dCell.Style = $dStyle
$ws.Range('E38').Value = '  -0.17%  '

$ws.Range('B39').Value = 'ThetaToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$dCell = $ws.Cells.Item(39, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '3.21'
$dCell.Style = $dStyle
$ws.Range('E39').Value = '  +12.83%  '

$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$dCell = $ws.Cells.Item(40, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.0₃0760'
$dCell.Style = $dStyle
$ws.Range('E40').Value = '  +3.41%  '

$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$dCell = $ws.Cells.Item(41, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '3.190.22'
$dCell.Style = $dStyle
$ws.Range('E41').Value = '  +4.89%  '

$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$dCell = $ws.Cells.Item(42, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.133'
$dCell.Style = $dStyle
$ws.Range('E42').Value = '  +1.23%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$dCell = $ws.Cells.Item(43, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.996'
$dCell.Style = $dStyle
$ws.Range('E43').Value = '  -0.13%  '

$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$dCell = $ws.Cells.Item(44, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '2.69'
$dCell.Style = $dStyle
$ws.Range('E44').Value = '  +7.56%  '

$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$dCell = $ws.Cells.Item(45, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '2.81'
$dCell.Style = $dStyle
$ws.Range('E45').Value = '  +14.05%  '

$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$dCell = $ws.Cells.Item(46, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '2.79'
$dCell.Style = $dStyle
$ws.Range('E46').Value = '  +19.62%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$dCell = $ws.Cells.Item(47, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.0410'
$dCell.Style = $dStyle
$ws.Range('E47').Value = '  +1.80%  '

$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$dCell = $ws.Cells.Item(48, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '0.131'
$dCell.Style = $dStyle
$ws.Range('E48').Value = '  +1.34%  '

$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$dCell = $ws.Cells.Item(49, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '3.07'
$dCell.Style = $dStyle
$ws.Range('E49').Value = '  -1.44%  '

$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$dCell = $ws.Cells.Item(50, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '8.60'
$dCell.Style = $dStyle
$ws.Range('E50').Value = '  +3.53%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$dCell = $ws.Cells.Item(51, 4)
$dStyle = $dCell.Style
$dCell.NumberFormat = '@'
$dCell.Value = '138.55'
$dCell.Style = $dStyle
$ws.Range('E51').Value = '  +0.98%  '
